$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values (D2:D51) remain text, matching the source data
# which stores numeric-looking prices as text (avoids Excel auto-number conversion).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.243.67"
$ws.Range("E2").Value = "  +4.28%  "
$ws.Range("D3").Value = "1.785.08"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "337.67"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "0.3826"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "0.3432"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "47.39"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "1.158"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").Value = "0.07422"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "23.52"
$ws.Range("E12").Value = "  +8.83%  "
$ws.Range("D13").Value = "0.9995"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "6.438"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "1.785.38"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "7.153"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "0.00001082"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").Value = "0.06662"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "82.83"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "0.9967"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "6.448"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "28.247.93"
$ws.Range("E23").Value = "  +4.23%  "
$ws.Range("D24").Value = "12.12"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").Value = "2.364"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").Value = "20.89"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").Value = "1.434"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "2.415"
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").Value = "154.77"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "136.14"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("B31").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C31").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D31").Value = "1.987.47"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "6.144"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").Value = "3.966"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").Value = "0.08861"
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("D35").Value = "12.80"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("D36").Value = "0.02440"
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("D37").Value = "0.6875"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").Value = "5.328"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").Value = "0.06375"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").Value = "0.2180"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.241"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.504"
$ws.Range("E42").Value = "  -7.39%  "
$ws.Range("D43").Value = "8.332"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").Value = "14.23"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").Value = "0.9972"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").Value = "0.6315"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").Value = "3.863"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "133.02"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("D49").Value = "2.098"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").Value = "0.07467"
$ws.Range("E50").Value = "  +5.38%  "
$ws.Range("D51").Value = "1.202"
$ws.Range("E51").Value = "  +8.30%  "
